# Weekly data refresh: a new price record (week of 2023-01-20) is inserted
# at row 563, pushing the existing rows 563..625 down to 564..626.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 563 (shifts 563:625 -> 564:626,
# dimension grows from A1:R625 to A1:R626).
$ws.Rows.Item(563).Insert()

# Populate the new row 563 with the new weekly record.
$ws.Cells.Item(563, 1).Value  = 3
$ws.Cells.Item(563, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(563, 3).Value  = "Coquimbo"
$ws.Cells.Item(563, 4).Value  = 44946
$ws.Cells.Item(563, 5).Value  = 5
$ws.Cells.Item(563, 6).Value  = 100112037
$ws.Cells.Item(563, 7).Value  = "Cebollín"
$ws.Cells.Item(563, 8).Value  = "Sin especificar"
$ws.Cells.Item(563, 9).Value  = "Primera"
$ws.Cells.Item(563, 10).Value = 230
$ws.Cells.Item(563, 11).Value = 3000
$ws.Cells.Item(563, 12).Value = 3500
$ws.Cells.Item(563, 13).Value = 3261
$ws.Cells.Item(563, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(563, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(563, 16).Value = 91
$ws.Cells.Item(563, 17).Value = 36
$ws.Cells.Item(563, 18).Value = "Hortaliza"
